$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8279455900192261
$ws.Range("B1").Value = 1.051031470298767
$ws.Range("C1").Value = 1.550809502601624
$ws.Range("D1").Value = 2.265527963638306
$ws.Range("E1").Value = 1.694389820098877
